# Added full correlation analysis and summary statistics
$wb = $excel.ActiveWorkbook

# --- stationarity_results: PP row (row 3) was missing lags_used/n_obs ---
$ws1 = $wb.Worksheets.Item("stationarity_results")
$ws1.Range("E3").Value = 31
$ws1.Range("F3").Value = 4424

# --- arch_lm_results: row 2 was missing LM_stat/LM_p_value/F_stat/F_p_value/n_obs,
#     and the conclusion changes now that ARCH effects are detected ---
$ws3 = $wb.Worksheets.Item("arch_lm_results")
$ws3.Range("C2").Value = 718.6642156973709
$ws3.Range("D2").Value = [double]"4.450569524860884e-146"
$ws3.Range("E2").Value = 71.32817042295778
$ws3.Range("F2").Value = [double]"7.333424758267242e-160"
$ws3.Range("G2").Value = 4425
$ws3.Range("H2").Value = "Reject no-ARCH (vol clustering)"

# --- garch_results: record the fitted innovation distribution ---
$ws4 = $wb.Worksheets.Item("garch_results")
$ws4.Range("I2").Value = "Standardized Student's t"

# --- forecast_results: add out-of-sample error / hit-rate summary columns
#     ahead of the existing window / ar_order_p columns ---
$ws6 = $wb.Worksheets.Item("forecast_results")
$ws6.Range("B1:F2").EntireColumn.Insert()

$ws6.Range("B1").Value = "n_obs"
$ws6.Range("C1").Value = "mse"
$ws6.Range("D1").Value = "mae"
$ws6.Range("E1").Value = "directional_hits"
$ws6.Range("F1").Value = "directional_hit_rate"

# match the bold/bordered header formatting used by the rest of row 1
$ws6.Range("G1").Copy()
$ws6.Range("B1:F1").PasteSpecial(-4122)

$ws6.Range("B2").Value = 1409
$ws6.Range("C2").Value = 13.22533602837114
$ws6.Range("D2").Value = 2.009767936038296
$ws6.Range("E2").Value = 830
$ws6.Range("F2").Value = 0.5894886363636364
# new data cells should stay unstyled, like the pre-existing data row
$ws6.Range("B2:F2").Style = "Normal"
